$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update item-1 row (row 7) with the new product/sale figures ---

# Plain text cells (their style already formats numFmt 49 = Text), so a
# normal value write is stored as text automatically.
$ws.Range("C7").Value = "ETHOXA 250MG/5ML SYRUP 120ML"
$ws.Range("N7").Value = "99.00"
$ws.Range("Q7").Value = "2:0"

# L7 / P7 carry a numeric display format (165 / 2) but must keep holding
# literal text, exactly like the original report. Temporarily force the
# cell to Text so the write lands as a shared string, then restore the
# original number format so the style index is unchanged.
$fmtL7 = $ws.Range("L7").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $fmtL7

$fmtP7 = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "198.0000"
$ws.Range("P7").NumberFormat = $fmtP7

# --- Remove the old item-2 row (row 8) entirely; rows below shift up ---
$ws.Rows("8:8").Delete()

# Row 8 (new) keeps its original row height (it used to belong to the
# deleted item-2 row) even though its content now comes from the old
# totals row.
$ws.Rows(8).RowHeight = 24.75

# --- Update the totals figure that is now on row 8 (a real number) ---
$ws.Range("P8").Value = 198

# --- Update the footer timestamp that is now on row 9 ---
$ws.Range("A9").Value = "Wednesday, 18 June, 2025 12:31 AM"
